$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.40000333333333
$ws.Range("H2").Value = 58.20001
$ws.Range("I2").Value = 0.2041610870830937
$ws.Range("J2").Value = 0.2041610870830936
$ws.Range("M2").Value = 90.25004833333332
$ws.Range("N2").Value = 270.750145
$ws.Range("O2").Value = 0.8928575650827933
$ws.Range("P2").Value = 0.8928575650827932
$ws.Range("Q2").Value = 1750.851238500161
$ws.Range("R2").Value = 15757.66114650145
$ws.Range("S2").Value = 0.1822867710976671
$ws.Range("T2").Value = 0.1822867710976671

# Row 3
$ws.Range("G3").Value = 19.40000333333333
$ws.Range("H3").Value = 58.20001
$ws.Range("I3").Value = 0.2041610870830937
$ws.Range("J3").Value = 0.2041610870830936
$ws.Range("O3").Value = 0.05133510428912089
$ws.Range("P3").Value = 0.05133510428912089
$ws.Range("Q3").Value = 100.6656990298422
$ws.Range("R3").Value = 905.99129126858
$ws.Range("S3").Value = 0.0104806306971909
$ws.Range("T3").Value = 0.0104806306971909

# Row 4
$ws.Range("G4").Value = 19.40000333333333
$ws.Range("H4").Value = 58.20001
$ws.Range("I4").Value = 0.2041610870830937
$ws.Range("J4").Value = 0.2041610870830936
$ws.Range("M4").Value = 5.380476000000001
$ws.Range("N4").Value = 16.141428
$ws.Range("O4").Value = 0.05322987398968605
$ws.Range("P4").Value = 0.05322987398968604
$ws.Range("Q4").Value = 104.38125233492
$ws.Range("R4").Value = 939.43127101428
$ws.Range("S4").Value = 0.0108674689390304
$ws.Range("T4").Value = 0.01086746893903039

# Row 5
$ws.Range("G5").Value = 19.40000333333333
$ws.Range("H5").Value = 58.20001
$ws.Range("I5").Value = 0.2041610870830937
$ws.Range("J5").Value = 0.2041610870830936
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2605293333333333
$ws.Range("N5").Value = 0.781588
$ws.Range("O5").Value = 0.002577456638399696
$ws.Range("P5").Value = 0.002577456638399696
$ws.Range("Q5").Value = 5.054269935097778
$ws.Range("R5").Value = 45.48842941588
$ws.Range("S5").Value = 0.0005262163492052183
$ws.Range("T5").Value = 0.0005262163492052182

# Row 6
$ws.Range("I6").Value = 0.1299257642351539
$ws.Range("J6").Value = 0.1299257642351539
$ws.Range("M6").Value = 90.25004833333332
$ws.Range("N6").Value = 270.750145
$ws.Range("O6").Value = 0.8928575650827933
$ws.Range("P6").Value = 0.8928575650827932
$ws.Range("Q6").Value = 1114.221561387036
$ws.Range("R6").Value = 10027.99405248332
$ws.Range("S6").Value = 0.1160052014965206
$ws.Range("T6").Value = 0.1160052014965206

# Row 7
$ws.Range("I7").Value = 0.1299257642351539
$ws.Range("J7").Value = 0.1299257642351539
$ws.Range("O7").Value = 0.05133510428912089
$ws.Range("P7").Value = 0.05133510428912089
$ws.Range("S7").Value = 0.006669752656855359
$ws.Range("T7").Value = 0.006669752656855358

# Row 8
$ws.Range("I8").Value = 0.1299257642351539
$ws.Range("J8").Value = 0.1299257642351539
$ws.Range("M8").Value = 5.380476000000001
$ws.Range("N8").Value = 16.141428
$ws.Range("O8").Value = 0.05322987398968605
$ws.Range("P8").Value = 0.05322987398968604
$ws.Range("Q8").Value = 66.42702669347202
$ws.Range("R8").Value = 597.843240241248
$ws.Range("S8").Value = 0.006915932058250901
$ws.Range("T8").Value = 0.006915932058250899

# Row 9
$ws.Range("I9").Value = 0.1299257642351539
$ws.Range("J9").Value = 0.1299257642351539
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2605293333333333
$ws.Range("N9").Value = 0.781588
$ws.Range("O9").Value = 0.002577456638399696
$ws.Range("P9").Value = 0.002577456638399696
$ws.Range("Q9").Value = 3.216479170200889
$ws.Range("R9").Value = 28.948312531808
$ws.Range("S9").Value = 0.0003348780235270513
$ws.Range("T9").Value = 0.0003348780235270512

# Row 10
$ws.Range("G10").Value = 62.84232966666667
$ws.Range("H10").Value = 188.526989
$ws.Range("I10").Value = 0.6613379451093298
$ws.Range("J10").Value = 0.6613379451093298
$ws.Range("M10").Value = 90.25004833333332
$ws.Range("N10").Value = 270.750145
$ws.Range("O10").Value = 0.8928575650827933
$ws.Range("P10").Value = 0.8928575650827932
$ws.Range("Q10").Value = 5671.523289795934
$ws.Range("R10").Value = 51043.7096081634
$ws.Range("S10").Value = 0.5904805873671742
$ws.Range("T10").Value = 0.5904805873671741

# Row 11
$ws.Range("G11").Value = 62.84232966666667
$ws.Range("H11").Value = 188.526989
$ws.Range("I11").Value = 0.6613379451093298
$ws.Range("J11").Value = 0.6613379451093298
$ws.Range("O11").Value = 0.05133510428912089
$ws.Range("P11").Value = 0.05133510428912089
$ws.Range("Q11").Value = 326.0858741033958
$ws.Range("R11").Value = 2934.772866930562
$ws.Range("S11").Value = 0.03394985238254035
$ws.Range("T11").Value = 0.03394985238254035

# Row 12
$ws.Range("G12").Value = 62.84232966666667
$ws.Range("H12").Value = 188.526989
$ws.Range("I12").Value = 0.6613379451093298
$ws.Range("J12").Value = 0.6613379451093298
$ws.Range("M12").Value = 5.380476000000001
$ws.Range("N12").Value = 16.141428
$ws.Range("O12").Value = 0.05322987398968605
$ws.Range("P12").Value = 0.05322987398968604
$ws.Range("Q12").Value = 338.1216465555881
$ws.Range("R12").Value = 3043.094819000292
$ws.Range("S12").Value = 0.03520293548276753
$ws.Range("T12").Value = 0.03520293548276753

# Row 13
$ws.Range("G13").Value = 62.84232966666667
$ws.Range("H13").Value = 188.526989
$ws.Range("I13").Value = 0.6613379451093298
$ws.Range("J13").Value = 0.6613379451093298
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2605293333333333
$ws.Range("N13").Value = 0.781588
$ws.Range("O13").Value = 0.002577456638399696
$ws.Range("P13").Value = 0.002577456638399696
$ws.Range("Q13").Value = 16.37227025317022
$ws.Range("R13").Value = 147.350432278532
$ws.Range("S13").Value = 0.001704569876847656
$ws.Range("T13").Value = 0.001704569876847656

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.4347496666666666
$ws.Range("H14").Value = 1.304249
$ws.Range("I14").Value = 0.004575203572422716
$ws.Range("J14").Value = 0.004575203572422716
$ws.Range("M14").Value = 90.25004833333332
$ws.Range("N14").Value = 270.750145
$ws.Range("O14").Value = 0.8928575650827933
$ws.Range("P14").Value = 0.8928575650827932
$ws.Range("Q14").Value = 39.23617842956722
$ws.Range("R14").Value = 353.125605866105
$ws.Range("S14").Value = 0.004085005121431444
$ws.Range("T14").Value = 0.004085005121431443

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.4347496666666666
$ws.Range("H15").Value = 1.304249
$ws.Range("I15").Value = 0.004575203572422716
$ws.Range("J15").Value = 0.004575203572422716
$ws.Range("O15").Value = 0.05133510428912089
$ws.Range("P15").Value = 0.05133510428912089
$ws.Range("Q15").Value = 2.255895442182444
$ws.Range("R15").Value = 20.303058979642
$ws.Range("S15").Value = 0.0002348685525342786
$ws.Range("T15").Value = 0.0002348685525342786

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.4347496666666666
$ws.Range("H16").Value = 1.304249
$ws.Range("I16").Value = 0.004575203572422716
$ws.Range("J16").Value = 0.004575203572422716
$ws.Range("M16").Value = 5.380476000000001
$ws.Range("N16").Value = 16.141428
$ws.Range("O16").Value = 0.05322987398968605
$ws.Range("P16").Value = 0.05322987398968604
$ws.Range("Q16").Value = 2.339160147508
$ws.Range("R16").Value = 21.052441327572
$ws.Range("S16").Value = 0.0002435375096372226
$ws.Range("T16").Value = 0.0002435375096372226

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.4347496666666666
$ws.Range("H17").Value = 1.304249
$ws.Range("I17").Value = 0.004575203572422716
$ws.Range("J17").Value = 0.004575203572422716
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.2605293333333333
$ws.Range("N17").Value = 0.781588
$ws.Range("O17").Value = 0.002577456638399696
$ws.Range("P17").Value = 0.002577456638399696
$ws.Range("Q17").Value = 0.1132650408235556
$ws.Range("R17").Value = 1.019385367412
$ws.Range("S17").Value = 0.00001179238881977094
$ws.Range("T17").Value = 0.00001179238881977093
